# affichage graphique dans Excel
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix mislabeled header cells: the "Contenu du stage" / "Type entreprise"
# labels were entered one column too far left (C16/C25 instead of B16/B25,
# matching the B10 "Lieu du stage" label above them).
$c16 = $ws.Range("C16").Value()
$ws.Range("B16").Value = $c16
$ws.Range("C16").ClearContents()

$c25 = $ws.Range("C25").Value()
$ws.Range("B25").Value = $c25
$ws.Range("C25").ClearContents()

# --- Chart 1: "Lieu du stage" pie chart (rows 10-14)
$co1 = $ws.ChartObjects().Add(0, 0, 0, 0)
$ch1 = $co1.Chart
$ch1.ChartType = 5
$ch1.SetSourceData($ws.Range("D10:E14"))
$ch1.HasTitle = $true
$ch1.ChartTitle.Text = "Lieu du stage"
$s1 = $ch1.SeriesCollection(1)
$s1.Name = "=Worksheet!`$B`$10"
$s1.HasDataLabels = $true
$dl1 = $s1.DataLabels()
$dl1.ShowValue = $true
$dl1.ShowPercentage = $true
$dl1.ShowLegendKey = $false
$dl1.ShowCategoryName = $false
$dl1.ShowSeriesName = $false
$dl1.ShowBubbleSize = $false
$pt1 = $s1.Points().Item(4)
$pt1.Format.Fill.Solid()
$pt1.Format.Fill.ForeColor.RGB = 39423
$ch1.HasLegend = $true
$ch1.Legend.Position = -4152

# --- Chart 2: "Contenu du stage" pie chart (rows 16-23)
$co2 = $ws.ChartObjects().Add(0, 0, 0, 0)
$ch2 = $co2.Chart
$ch2.ChartType = 5
$ch2.SetSourceData($ws.Range("D16:E23"))
$ch2.HasTitle = $true
$ch2.ChartTitle.Text = "Contenu du stage"
$s2 = $ch2.SeriesCollection(1)
$s2.Name = "=Worksheet!`$B`$16"
$s2.HasDataLabels = $true
$dl2 = $s2.DataLabels()
$dl2.ShowValue = $true
$dl2.ShowPercentage = $true
$dl2.ShowLegendKey = $false
$dl2.ShowCategoryName = $false
$dl2.ShowSeriesName = $false
$dl2.ShowBubbleSize = $false
$pt2 = $s2.Points().Item(4)
$pt2.Format.Fill.Solid()
$pt2.Format.Fill.ForeColor.RGB = 39423
$ch2.HasLegend = $true
$ch2.Legend.Position = -4152

# --- Chart 3: "Type du stage" pie chart (rows 25-28)
$co3 = $ws.ChartObjects().Add(0, 0, 0, 0)
$ch3 = $co3.Chart
$ch3.ChartType = 5
$ch3.SetSourceData($ws.Range("D25:E28"))
$ch3.HasTitle = $true
$ch3.ChartTitle.Text = "Type du stage"
$s3 = $ch3.SeriesCollection(1)
$s3.Name = "=Worksheet!`$B`$25"
$s3.HasDataLabels = $true
$dl3 = $s3.DataLabels()
$dl3.ShowValue = $true
$dl3.ShowPercentage = $true
$dl3.ShowLegendKey = $false
$dl3.ShowCategoryName = $false
$dl3.ShowSeriesName = $false
$dl3.ShowBubbleSize = $false
$pt3 = $s3.Points().Item(4)
$pt3.Format.Fill.Solid()
$pt3.Format.Fill.ForeColor.RGB = 39423
$ch3.HasLegend = $true
$ch3.Legend.Position = -4152

Write-Host "done"
